$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MSE")

# --- Row 3: header titles ---
$ws.Range("D3").Copy()
$ws.Range("N3").PasteSpecial(-4122)
$ws.Range("F3").Copy()
$ws.Range("P3").PasteSpecial(-4122)
$ws.Range("N3").Value = "Random Forest-100 (superdataset-20.csv without cons)"
$ws.Range("D3").Value = "Random Forest-100 (superdataset-21.csv without cons)"
$ws.Range("I3").Value = "Random Forest-100 (superdataset-21.csv)"

# --- Row 4: sub headers ---
$ws.Range("D4").Copy()
$ws.Range("N4").PasteSpecial(-4122)
$ws.Range("E4").Copy()
$ws.Range("O4").PasteSpecial(-4122)
$ws.Range("F4").Copy()
$ws.Range("P4").PasteSpecial(-4122)
$ws.Range("O4").Value = "train (MSE)"
$ws.Range("P4").Value = "test (MSE)"

# --- Rows 5-54: data ---
$ws.Range("J5").Value = 0.0001150940191894486
$ws.Range("K5").Value = 0.0008459625171393564
$ws.Range("J6").Value = 0.0001089138120210779
$ws.Range("K6").Value = 0.0007794426216960208
$ws.Range("J7").Value = 0.0001015404112154605
$ws.Range("K7").Value = 0.001068389897890192
$ws.Range("J8").Value = 0.0001118512864312136
$ws.Range("K8").Value = 0.0009344658866757057
$ws.Range("J9").Value = 0.0001188819329491765
$ws.Range("K9").Value = 0.0007126229027877369
$ws.Range("J10").Value = 0.0001103963584647762
$ws.Range("K10").Value = 0.0007885029855490761
$ws.Range("J11").Value = 0.0001170134769518971
$ws.Range("K11").Value = 0.0008236579425686637
$ws.Range("J12").Value = 0.0001170008588759945
$ws.Range("K12").Value = 0.0008130380797629764
$ws.Range("J13").Value = 0.0001074416458598002
$ws.Range("K13").Value = 0.0008674383149265965
$ws.Range("J14").Value = 0.000109577897553583
$ws.Range("K14").Value = 0.0007080205048467555
$ws.Range("J15").Value = 0.000111465664733632
$ws.Range("K15").Value = 0.0008952242731027751
$ws.Range("J16").Value = 0.0001288116014403564
$ws.Range("K16").Value = 0.0006054788756640541
$ws.Range("J17").Value = 0.0001160116806879669
$ws.Range("K17").Value = 0.0007868620129807003
$ws.Range("J18").Value = 0.0001151920989721826
$ws.Range("K18").Value = 0.0006752920261545541
$ws.Range("J19").Value = 0.0001113182466093625
$ws.Range("K19").Value = 0.0006464739280765895
$ws.Range("J20").Value = 0.0001088901512947554
$ws.Range("K20").Value = 0.0008073235333991647
$ws.Range("J21").Value = 0.0001101115699044487
$ws.Range("K21").Value = 0.0009227049275820674
$ws.Range("J22").Value = 0.0001043096559125857
$ws.Range("K22").Value = 0.0009052808885131904
$ws.Range("J23").Value = 0.0001194774964722091
$ws.Range("K23").Value = 0.0006369086269405749
$ws.Range("J24").Value = 0.0001158794379207492
$ws.Range("K24").Value = 0.0006649978987220839
$ws.Range("J25").Value = 0.0001200670400042585
$ws.Range("K25").Value = 0.0008228131873474324
$ws.Range("J26").Value = 0.0001115945781844151
$ws.Range("K26").Value = 0.0009152913625696437
$ws.Range("J27").Value = 0.0001129101404089035
$ws.Range("K27").Value = 0.0008687169280022576
$ws.Range("J28").Value = 0.0001032695086839765
$ws.Range("K28").Value = 0.0009360304800874892
$ws.Range("J29").Value = 0.0001146587046893823
$ws.Range("K29").Value = 0.0007971683732901313
$ws.Range("J30").Value = 0.000102154612558613
$ws.Range("K30").Value = 0.001333007131739747
$ws.Range("J31").Value = 0.000112590788315525
$ws.Range("K31").Value = 0.001020304846375869
$ws.Range("J32").Value = 0.0001084165222681405
$ws.Range("K32").Value = 0.000786724981854342
$ws.Range("J33").Value = 0.0001222480652536344
$ws.Range("K33").Value = 0.0006188467374870361
$ws.Range("J34").Value = 0.0001107273227631967
$ws.Range("K34").Value = 0.000777346114224395
$ws.Range("J35").Value = 0.0001178908565366489
$ws.Range("K35").Value = 0.0006614372410414719
$ws.Range("J36").Value = 0.0001061911727931064
$ws.Range("K36").Value = 0.001083629364751326
$ws.Range("J37").Value = 0.0001158819714449737
$ws.Range("K37").Value = 0.0007756700066746047
$ws.Range("J38").Value = 0.0001214944951135354
$ws.Range("K38").Value = 0.0007118090664459475
$ws.Range("J39").Value = 0.0001100428270499541
$ws.Range("K39").Value = 0.0008365367047279698
$ws.Range("J40").Value = 0.0001174069389565614
$ws.Range("K40").Value = 0.0008172358225092515
$ws.Range("J41").Value = 0.000114375185014424
$ws.Range("K41").Value = 0.0008602366107644186
$ws.Range("J42").Value = 0.0001159353146501999
$ws.Range("K42").Value = 0.0007270020467061802
$ws.Range("J43").Value = 0.0001111612736735101
$ws.Range("K43").Value = 0.0007613728967302548
$ws.Range("J44").Value = 0.0001198889625765995
$ws.Range("K44").Value = 0.0006617120280276574
$ws.Range("J45").Value = 0.0001104471776958471
$ws.Range("K45").Value = 0.0007302408854739465
$ws.Range("J46").Value = 0.0001059238435859602
$ws.Range("K46").Value = 0.000971503809600807
$ws.Range("J47").Value = 0.000122426004213166
$ws.Range("K47").Value = 0.000814928426148509
$ws.Range("J48").Value = 0.0001161235099756959
$ws.Range("K48").Value = 0.0007057338061830683
$ws.Range("J49").Value = 0.000119698277375593
$ws.Range("K49").Value = 0.0007341959865992744
$ws.Range("J50").Value = 0.0001190421091684601
$ws.Range("K50").Value = 0.0006527722915068563
$ws.Range("J51").Value = 0.0001135124577623107
$ws.Range("K51").Value = 0.001007523653543452
$ws.Range("J52").Value = 0.0001113164752526531
$ws.Range("K52").Value = 0.0008714348231375634
$ws.Range("J53").Value = 0.0001162135403229242
$ws.Range("K53").Value = 0.0006580047869099962
$ws.Range("J54").Value = 0.000104656791805727
$ws.Range("K54").Value = 0.001035666480828549

# --- N column (counter like D/I), O/P columns (style like J/K, empty) ---
$ws.Range("D5").Copy()
$ws.Range("N5").PasteSpecial(-4122)
$ws.Range("N5").Value = 1
$ws.Range("J5").Copy()
$ws.Range("O5").PasteSpecial(-4122)
$ws.Range("K5").Copy()
$ws.Range("P5").PasteSpecial(-4122)

$ws.Range("D6").Copy()
$ws.Range("N6").PasteSpecial(-4122)
$ws.Range("N6").Formula = "=N5+1"
$ws.Range("J6").Copy()
$ws.Range("O6").PasteSpecial(-4122)
$ws.Range("K6").Copy()
$ws.Range("P6").PasteSpecial(-4122)

$ws.Range("D7").Copy()
$ws.Range("N7:N54").PasteSpecial(-4122)
$ws.Range("N7:N54").Formula = "=N6+1"

for ($r = 7; $r -le 54; $r++) {
    $srcJ = "J" + $r
    $srcK = "K" + $r
    $dstO = "O" + $r
    $dstP = "P" + $r
    $ws.Range($srcJ).Copy()
    $ws.Range($dstO).PasteSpecial(-4122)
    $ws.Range($srcK).Copy()
    $ws.Range($dstP).PasteSpecial(-4122)
}

# --- Row 56: avg ---
$ws.Range("D56").Copy()
$ws.Range("N56").PasteSpecial(-4122)
$ws.Range("N56").Value = "avg"
$ws.Range("J56").Copy()
$ws.Range("O56").PasteSpecial(-4122)
$ws.Range("O56").Formula = "=AVERAGE(O5:O54)"
$ws.Range("K56").Copy()
$ws.Range("P56").PasteSpecial(-4122)
$ws.Range("P56").Formula = "=AVERAGE(P5:P54)"

# --- Row 57: SD ---
$ws.Range("D57").Copy()
$ws.Range("N57").PasteSpecial(-4122)
$ws.Range("N57").Value = "SD"
$ws.Range("J57").Copy()
$ws.Range("O57").PasteSpecial(-4122)
$ws.Range("O57").Formula = "=_xlfn.STDEV.S(O5:O54)"
$ws.Range("K57").Copy()
$ws.Range("P57").PasteSpecial(-4122)
$ws.Range("P57").Formula = "=_xlfn.STDEV.S(P5:P54)"

# --- Column widths ---
$ws.Columns.Item(15).ColumnWidth = 11.1
$ws.Columns.Item(16).ColumnWidth = 13.65

# --- Recalculate ---
$wb.Application.Calculate()

# --- View / selection ---
$ws.Activate()
$ws.Range("O17").Select()

Write-Output "done"
